# Applies the content edits described by the commit diff to DaC.docx.
#
# Three passages of "VISTO ..." legal-citation text are updated with new
# provvedimento numbers / protocol numbers / dates (PNRR boilerplate refresh).
# A cosmetic styles.xml flag (<w:semiHidden/> on the "Default Paragraph Font"
# style) is not reachable through the Word object model exposed by this host
# (Style.Hidden / Style.Visibility only round-trip <w:hidden>, never
# <w:semiHidden>), so it is intentionally left untouched.

$d = $word.ActiveDocument

$wdFindContinue = 1
$wdReplaceOne = 1

function Replace-Text($find, $replace) {
    $rng = $d.Content
    $ok = $rng.Find.Execute(
        $find, $false, $false, $false, $false, $false, $true,
        $wdFindContinue, $false, $replace, $wdReplaceOne)
    Write-Output "Replace [$find] -> [$replace] : $ok"
}

# --- Change 1 ----------------------------------------------------------
# "Regolamento di Organizzazione e Funzionamento del CNR emanato con
#  Provvedimento del Presidente n. 119, prot. n. 241776 del 10/07/2024,
#  in vigore dal 1° agosto 2024"
#    -> "... Provvedimento del Presidente n. 144 Prot. n. 521963 del
#  19 dicembre 2025, in vigore dal 1° gennaio 2026"
Replace-Text `
    "Provvedimento del Presidente n. 119, prot. n. 241776 del 10/07/2024, in vigore dal 1° agosto 2024" `
    "Provvedimento del Presidente n. 144 Prot. n. 521963 del 19 dicembre 2025, in vigore dal 1° gennaio 2026"

# --- Change 2 ------------------------------------------------------------
# "... entrato in vigore dal 1° gennaio 2025 ed," -> "... 1° gennaio 2025 e,"
# (also collapses the two proofErr-wrapped "ed" / "," runs into plain text)
Replace-Text `
    "entrato in vigore dal 1° gennaio 2025 ed," `
    "entrato in vigore dal 1° gennaio 2025 e,"

# --- Change 3 ------------------------------------------------------------
# "il bilancio di previsione del Consiglio Nazionale delle Ricerche per
#  l'esercizio finanziario 2025, approvato dal Consiglio di Amministrazione
#  con deliberazione n° 420/2024 del 17/12/2024"
#    -> "il Bilancio Unico di Previsione del Consiglio Nazionale delle
#  Ricerche per l'esercizio finanziario 2026, approvato dal Consiglio di
#  Amministrazione con deliberazione n. 245/2025 – Verbale 527 del
#  17 dicembre 2025"
Replace-Text `
    "il bilancio di previsione del Consiglio Nazionale delle Ricerche per l'esercizio finanziario 2025, approvato dal Consiglio di Amministrazione con deliberazione n° 420/2024 del 17/12/2024" `
    "il Bilancio Unico di Previsione del Consiglio Nazionale delle Ricerche per l’esercizio finanziario 2026, approvato dal Consiglio di Amministrazione con deliberazione n. 245/2025 – Verbale 527 del 17 dicembre 2025"
